# Natmi following Dr Hou advice:
# recompute the LR-pair stats for the existing Sending-cluster groups
# (M1, M2, Neutro) against Cd28/Cd86 and add a new "sCs" sending-cluster
# group (rows 11-13) with its own stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("M1", "Cd28", "Cd86", "M1", 3, 1, 4.430875666666666, 13.292627, 0.3569860433688579, 0.3569860433688579, 3, 1, 81.47348266666667, 244.420448, 0.3594530042390097, 0.3594530042390096, 360.9988718263218, 3248.989846436896, 0.1283197057603334, 0.1283197057603333),
  @("M1", "Cd28", "Cd86", "M2", 3, 1, 4.430875666666666, 13.292627, 0.3569860433688579, 0.3569860433688579, 3, 1, 123.2241746666667, 369.672524, 0.5436529571225457, 0.5436529571225457, 545.9909970756164, 4913.918973680548, 0.194076518128957, 0.1940765181289569),
  @("M1", "Cd28", "Cd86", "Neutro", 3, 1, 4.430875666666666, 13.292627, 0.3569860433688579, 0.3569860433688579, 3, 1, 21.96196633333333, 65.885899, 0.0968940386384447, 0.09689403863844469, 97.31074221851921, 875.7966799666729, 0.03458981947956762, 0.0345898194795676),
  @("M2", "Cd28", "Cd86", "M1", 3, 1, 7.025118333333332, 21.075355, 0.5659985489733577, 0.5659985489733576, 3, 1, 81.47348266666667, 244.420448, 0.3594530042390097, 0.3594530042390096, 572.3608567621155, 5151.247710859039, 0.2034498788233937, 0.2034498788233936),
  @("M2", "Cd28", "Cd86", "M2", 3, 1, 7.025118333333332, 21.075355, 0.5659985489733577, 0.5659985489733576, 3, 1, 123.2241746666667, 369.672524, 0.5436529571225457, 0.5436529571225457, 865.6644085606688, 7790.97967704602, 0.3077067848764359, 0.3077067848764358),
  @("M2", "Cd28", "Cd86", "Neutro", 3, 1, 7.025118333333332, 21.075355, 0.5659985489733577, 0.5659985489733576, 3, 1, 21.96196633333333, 65.885899, 0.0968940386384447, 0.09689403863844469, 154.2854123243494, 1388.568710919145, 0.05484188527352815, 0.05484188527352814),
  @("Neutro", "Cd28", "Cd86", "M1", 3, 1, 0.950574, 2.851722, 0.0765856857013987, 0.07658568570139869, 3, 1, 81.47348266666667, 244.420448, 0.3594530042390097, 0.3594530042390096, 77.446574312384, 697.019168811456, 0.02752895480707233, 0.02752895480707232),
  @("Neutro", "Cd28", "Cd86", "M2", 3, 1, 0.950574, 2.851722, 0.0765856857013987, 0.07658568570139869, 3, 1, 123.2241746666667, 369.672524, 0.5436529571225457, 0.5436529571225457, 117.133696609592, 1054.203269486328, 0.04163603450482327, 0.04163603450482326),
  @("Neutro", "Cd28", "Cd86", "Neutro", 3, 1, 0.950574, 2.851722, 0.0765856857013987, 0.07658568570139869, 3, 1, 21.96196633333333, 65.885899, 0.0968940386384447, 0.09689403863844469, 20.876474185342, 187.888267668078, 0.007420696389503108, 0.007420696389503105),
  @("sCs", "Cd28", "Cd86", "M1", 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.0004297219563856788, 0.0004297219563856788, 3, 1, 81.47348266666667, 244.420448, 0.3594530042390097, 0.3594530042390096, 0.4345523987164445, 3.910971588448, 0.000154464848210297, 0.0001544648482102969),
  @("sCs", "Cd28", "Cd86", "M2", 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.0004297219563856788, 0.0004297219563856788, 3, 1, 123.2241746666667, 369.672524, 0.5436529571225457, 0.5436529571225457, 0.6572366729471112, 5.915130056524001, 0.0002336196123295599, 0.0002336196123295599),
  @("sCs", "Cd28", "Cd86", "Neutro", 1, 0.3333333333333333, 0.005333666666666667, 0.016001, 0.0004297219563856788, 0.0004297219563856788, 3, 1, 21.96196633333333, 65.885899, 0.0968940386384447, 0.09689403863844469, 0.1171378077665556, 1.054240269899, 0.00004163749584582201, 0.000041637495845822)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

